$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3148.3125
$ws.Range("I40").Value = 2490.3076
$ws.Range("J40").Value = 5999.6665
$ws.Range("K40").Value = 2490.3076
$ws.Range("L40").Value = 5999.6665
$ws.Range("M40").Value = -2315.3076
$ws.Range("N40").Value = -6349.6665
$ws.Range("H69").Value = 32142.5
$ws.Range("J69").Value = 39856.668
$ws.Range("L69").Value = 119570.004
$ws.Range("N69").Value = -121318.004
$ws.Range("H72").Value = 32142.5
$ws.Range("J72").Value = 39856.668
$ws.Range("L72").Value = 358710.012
$ws.Range("N72").Value = -367446.012
$ws.Range("H93").Value = 200350000
$ws.Range("J93").Value = 200350000
$ws.Range("L93").Value = 200350000
$ws.Range("N93").Value = -200354992
$ws.Range("H100").Value = 100090.93
$ws.Range("I100").Value = 167325.83
$ws.Range("J100").Value = 55267.668
$ws.Range("K100").Value = 167325.83
$ws.Range("L100").Value = 55267.668
$ws.Range("M100").Value = -166784.83
$ws.Range("N100").Value = -56349.668
$ws.Range("H112").Value = 3834.6667
$ws.Range("I112").Value = 1103.2
$ws.Range("J112").Value = 4885.231
$ws.Range("K112").Value = 3309.6
$ws.Range("L112").Value = 14655.693
$ws.Range("M112").Value = -2201.6
$ws.Range("N112").Value = -16871.693
$ws.Range("H134").Value = 88652.39999999999
$ws.Range("J134").Value = 88652.39999999999
$ws.Range("L134").Value = 88652.39999999999
$ws.Range("N134").Value = -98792.39999999999
$ws.Range("H141").Value = 7095.875
$ws.Range("J141").Value = 10000
$ws.Range("L141").Value = 30000
$ws.Range("N141").Value = -40360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2918.5862
$ws.Range("I32").Value = 2918.5862
$ws.Range("K32").Value = 2918.5862
$ws.Range("M32").Value = -2631.5862
$ws.Range("H74").Value = 5504.4375
$ws.Range("I74").Value = 1922.5
$ws.Range("K74").Value = 1922.5
$ws.Range("M74").Value = -1048.5
$ws.Range("H77").Value = 5504.4375
$ws.Range("I77").Value = 1922.5
$ws.Range("K77").Value = 9612.5
$ws.Range("M77").Value = -5244.5
$ws.Range("H104").Value = 22222
$ws.Range("J104").Value = 22222
$ws.Range("L104").Value = 22222
$ws.Range("N104").Value = -29210
$ws.Range("H122").Value = 472278.03
$ws.Range("I122").Value = 4768.8887
$ws.Range("J122").Value = 1173541.8
$ws.Range("K122").Value = 14306.6661
$ws.Range("L122").Value = 3520625.4
$ws.Range("M122").Value = -11856.6661
$ws.Range("N122").Value = -3525525.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4192.1333
$ws.Range("I86").Value = 6383.125
$ws.Range("K86").Value = 6383.125
$ws.Range("M86").Value = -5260.125
$ws.Range("H89").Value = 4192.1333
$ws.Range("I89").Value = 6383.125
$ws.Range("K89").Value = 31915.625
$ws.Range("M89").Value = -26299.625
$ws.Range("H92").Value = 50167
$ws.Range("J92").Value = 50167
$ws.Range("L92").Value = 50167
$ws.Range("N92").Value = -55159
$ws.Range("H105").Value = 130749.875
$ws.Range("I105").Value = 202399.8
$ws.Range("J105").Value = 11333.333
$ws.Range("K105").Value = 202399.8
$ws.Range("L105").Value = 11333.333
$ws.Range("M105").Value = -200652.8
$ws.Range("N105").Value = -14827.333
$ws.Range("H134").Value = 6668.3447
$ws.Range("I134").Value = 6754.8887
$ws.Range("K134").Value = 20264.6661
$ws.Range("M134").Value = -17729.6661
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1287.375
$ws.Range("I22").Value = 996.5
$ws.Range("J22").Value = 1384.3334
$ws.Range("K22").Value = 996.5
$ws.Range("L22").Value = 1384.3334
$ws.Range("M22").Value = -646.5
$ws.Range("N22").Value = -2084.3334
$ws.Range("H99").Value = 7746791
$ws.Range("I99").Value = 23225422
$ws.Range("J99").Value = 7475.1
$ws.Range("K99").Value = 23225422
$ws.Range("L99").Value = 7475.1
$ws.Range("M99").Value = -23223924
$ws.Range("N99").Value = -10471.1
$ws.Range("I105").Value = 424640
$ws.Range("J105").Value = 4600
$ws.Range("K105").Value = 424640
$ws.Range("L105").Value = 4600
$ws.Range("M105").Value = -422893
$ws.Range("N105").Value = -8094
$ws.Range("H126").Value = 7746791
$ws.Range("I126").Value = 23225422
$ws.Range("J126").Value = 7475.1
$ws.Range("K126").Value = 69676266
$ws.Range("L126").Value = 22425.3
$ws.Range("M126").Value = -69673796
$ws.Range("N126").Value = -27365.3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 123.8125
$ws.Range("I12").Value = 281.16666
$ws.Range("K12").Value = 843.4999799999999
$ws.Range("M12").Value = -670.4999799999999
$ws.Range("H23").Value = 5747546.5
$ws.Range("J23").Value = 8772361
$ws.Range("L23").Value = 26317083
$ws.Range("N23").Value = -26317553
$ws.Range("H121").Value = 2042768.9
$ws.Range("I121").Value = 1808900.8
$ws.Range("J121").Value = 2224666.2
$ws.Range("K121").Value = 5426702.4
$ws.Range("L121").Value = 6673998.600000001
$ws.Range("M121").Value = -5425392.4
$ws.Range("N121").Value = -6676618.600000001
$ws.Range("H127").Value = 100
$ws.Range("J127").Value = 100
$ws.Range("L127").Value = 300
$ws.Range("N127").Value = -10220
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 21848.625
$ws.Range("I92").Value = 14999
$ws.Range("K92").Value = 14999
$ws.Range("M92").Value = -13127
$ws.Range("H105").Value = 55788.89
$ws.Range("J105").Value = 55788.89
$ws.Range("L105").Value = 55788.89
$ws.Range("N105").Value = -62776.89
$ws.Range("H122").Value = 8558.679
$ws.Range("I122").Value = 5165.864
$ws.Range("J122").Value = 20999
$ws.Range("K122").Value = 15497.592
$ws.Range("L122").Value = 62997
$ws.Range("M122").Value = -13047.592
$ws.Range("N122").Value = -67897
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 160000000
$ws.Range("I2").Value = 800000000
$ws.Range("K2").Value = 800000000
$ws.Range("M2").Value = -799999888
$ws.Range("H19").Value = 1766.6666
$ws.Range("I19").Value = 650
$ws.Range("K19").Value = 650
$ws.Range("M19").Value = -480
$ws.Range("H21").Value = 1500
$ws.Range("J21").Value = 1500
$ws.Range("L21").Value = 1500
$ws.Range("N21").Value = -1848
$ws.Range("H127").Value = 85000
$ws.Range("J127").Value = 85000
$ws.Range("L127").Value = 85000
$ws.Range("N127").Value = -94920
$ws.Range("H136").Value = 5321.9165
$ws.Range("I136").Value = 2719.2354
$ws.Range("J136").Value = 11642.714
$ws.Range("K136").Value = 8157.706200000001
$ws.Range("L136").Value = 34928.142
$ws.Range("M136").Value = -5607.706200000001
$ws.Range("N136").Value = -40028.142
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7216.522
$ws.Range("J122").Value = 12000
$ws.Range("L122").Value = 36000
$ws.Range("N122").Value = -40900
$ws.Range("H132").Value = 8625.138000000001
$ws.Range("I132").Value = 9435.106
$ws.Range("J132").Value = 4806.7144
$ws.Range("K132").Value = 28305.318
$ws.Range("L132").Value = 14420.1432
$ws.Range("M132").Value = -25775.318
$ws.Range("N132").Value = -19480.1432
